# Correccion a Diebold Mariano y revision de Cap1
#
# A "d=6" row was missing between "d=5" and what used to be labelled "d=7"
# (the old "d=7" row actually held the values that belong to "d=6"). The
# fix inserts the real "d=6" row (with freshly computed values) above the
# old row, pushing the previous "d=7" and "d=10" rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 8 ("d=10") down to row 9, and row 7 ("d=7") down to row 8,
# copying formatting and values separately so both travel with the row
# (PasteSpecial(xlPasteAll) alone does not carry the cell style here).
$ws.Range("A8:E8").Copy()
$ws.Range("A9:E9").PasteSpecial(-4122)
$ws.Range("A8:E8").Copy()
$ws.Range("A9:E9").PasteSpecial(-4163)

$ws.Range("A7:E7").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122)
$ws.Range("A7:E7").Copy()
$ws.Range("A8:E8").PasteSpecial(-4163)

# Row 7 becomes the new "d=6" row with the corrected values.
$ws.Range("A7").Value = "d=6"
$ws.Range("B7").Value = 97.85838074038192
$ws.Range("C7").Value = 97.94477153128256
$ws.Range("D7").Value = 97.97749483048452
$ws.Range("E7").Value = 97.9362554573154
